$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.740.92'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  -1.82%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.745.57'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  -4.36%  '

$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.60'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  -2.75%  '

$ws.Range('E6').Value = '  +0.08%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4201'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  -5.26%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3622'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').Value = '  -3.55%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.49'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  -5.15%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07408'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').Value = '  -3.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.084'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  -3.61%  '

$ws.Range('E12').Value = '  -0.08%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.60'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  -6.96%  '

$ws.Range('E14').Value = '  -4.69%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.259'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  -3.68%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.763.49'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  -3.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.40'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').Value = '  -3.26%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001049'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  -2.92%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06336'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  -2.33%  '

$ws.Range('E20').Value = '  +0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.95'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  -3.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.912'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  -6.44%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.774.15'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  -1.96%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.10'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  -5.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.095'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  -4.06%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.82'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  +0.83%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.06'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  -3.12%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.977.44'
$ws.Range('D28').Style = 'Normal'

$ws.Range('E28').Value = '  -2.95%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.123'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  -9.70%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.39'
$ws.Range('D30').Style = 'Normal'

$ws.Range('E30').Value = '  -3.84%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.118'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  -6.35%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.634'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  -0.89%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.536'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -5.77%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08785'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  -5.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.24'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  -6.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02264'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  -3.37%  '

$ws.Range('E37').Value = '  -4.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05990'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  -3.09%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.924'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  -4.84%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6276'
$ws.Range('D40').Style = 'Normal'

$ws.Range('E40').Value = '  -4.59%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.171'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  -2.96%  '

$ws.Range('E42').Value = '  +0.03%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.394'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  +0.17%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.803'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -3.85%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.30'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  -4.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5843'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  -3.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.669'
$ws.Range('D47').Style = 'Normal'

$ws.Range('E47').Value = '  -2.68%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.27'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  -3.76%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.958'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  -4.26%  '

$ws.Range('E50').Value = '  +1.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06787'
$ws.Range('D51').Style = 'Normal'

$ws.Range('E51').Value = '  -3.05%  '
